$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the changed values in rows 2-4
$ws.Range("B3").Value = "-1.871***"

$ws.Range("C2").Value = "0.046**"
# "-0.024" parses as a pure number, so force it to text with a quote
# prefix and then strip the resulting style back to Normal so the cell
# keeps the default (unstyled) look of its neighbours.
$ws.Range("C3").Value = "'-0.024"
$ws.Range("C3").Style = "Normal"
$ws.Range("C4").Value = "-0.024**"

$ws.Range("D2").Value = "0.33***"
$ws.Range("D3").Value = "7.831***"

# Remove the "Constant" and "r2_adj" rows (rows 5 and 6) entirely
$ws.Range("A5:D6").EntireRow.Delete()
